$wb = $excel.ActiveWorkbook

# Add a new worksheet at the very end of the workbook and name it "Android"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Android"

# Populate the new sheet with the Android best-practices content
$ws.Range("B2").Value = "Context"
$ws.Range("C3").Value = "Like scope that object lives in."
$ws.Range("B4").Value = "Activity/Fragment life cycle"
$ws.Range("C5").Value = "onCreate"
$ws.Range("D6").Value = "Create and init member data only"
$ws.Range("C7").Value = "onCreateView(Fragment)"
$ws.Range("D8").Value = "bind view only"
$ws.Range("C9").Value = "onstart"
$ws.Range("D10").Value = "set adapter for listview"
$ws.Range("D11").Value = "do some logic update stuff"
$ws.Range("C12").Value = "onpause"
$ws.Range("D13").Value = "save data"
$ws.Range("B14").Value = "Setting screen"
$ws.Range("C15").Value = "https://developer.android.com/guide/topics/ui/settings.html"
$ws.Range("B16").Value = "Divide layout(drawable,etc) into sub folders"
$ws.Range("C17").Value = "http://stackoverflow.com/questions/16577782/sub-folders-in-drawable-resource-folder"

# Match the final selection left on the new sheet (cell below the last entry)
$ws.Range("B18").Select() | Out-Null
